$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.651.88"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.473.64"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.79"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0870"
$ws.Range("E10").Value = "  +11.03%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "2.856.03"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.66"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "2.492.11"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("D18").Value = "41.612.70"
$ws.Range("D19").Value = "0.0₃0959"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.50"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.44"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.70"
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.84"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.91"
$ws.Range("E29").Value = "  +2.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.53"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.93"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0773"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.46"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.01"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.984.73"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.32"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").Value = "2.714.07"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.64"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.69"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.95"
$ws.Range("E51").Value = "  +2.44%  "
